# Add a new "Dansk titel" column before the existing "Status" column, and a
# new "Stamdata" column at the end; then populate the newly-added cells with
# their Danish-language / status data (per the Treasury Project Manager Role
# Center page-extension work: Table list.xlsx).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D. This shifts the existing "Status" column (D) to
# E and the existing "Comments in ReadMe" column (E) to F, exactly like
# selecting column D in Excel and choosing Insert.
$ws.Columns("D").Insert()

# New column headers.
$ws.Range("D1").Value = "Dansk titel"
$ws.Range("G1").Value = "Stamdata"

# Row 2 - Account / Depot
$ws.Range("D2").Value = "Depot"
$ws.Range("G2").Value = "Yes"

# Row 3 - Security / Værdipapir
$ws.Range("D3").Value = "Værdipapir"
$ws.Range("E3").Value = "Completed for now"
$ws.Range("F3").Value = "Yes"

# Row 4 - Security Ledger Entry / Værdipapirspost
$ws.Range("D4").Value = "Værdipapirspost"
$ws.Range("E4").Value = "Completed"

# Row 6 - Security Account / Depotbeholdning
$ws.Range("D6").Value = "Depotbeholdning"
$ws.Range("G6").Value = "Yes"

# Row 7 - Security Account Ledger Entry / Depotbeholdsningspost
$ws.Range("D7").Value = "Depotbeholdsningspost"
$ws.Range("E7").Value = "Completed"

# Row 10 - Security Register / Værdipapirjournal
$ws.Range("D10").Value = "Værdipapirjournal"
$ws.Range("E10").Value = "Security Register"

# Row 11 - Security Setup / Opsætning af Treasury
$ws.Range("D11").Value = "Opsætning af Treasury"
$ws.Range("G11").Value = "Yes"

# Row 12 - Security Posting Group / Værdipapirbogføringsgruppe
$ws.Range("D12").Value = "Værdipapirbogføringsgruppe"
$ws.Range("G12").Value = "Yes"

# Row 13 - Security Value Posting Group
$ws.Range("G13").Value = "Yes"

# Row 15 - Security Comment Line / Værdipapirbemærkningslinje
$ws.Range("D15").Value = "Værdipapirbemærkningslinje"
$ws.Range("E15").Value = "Completed"
$ws.Range("G15").Value = "Yes"

# Row 19 - Price / Fondskurser
$ws.Range("D19").Value = "Fondskurser"
$ws.Range("E19").Value = "Completed"
$ws.Range("F19").Value = "Yes"
